$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update of 2025 data and RF changes:
# The RF (raising factor) for the "2-RAP" sub-group rows 18:52 was recalculated.
$newRF = 20.64918032786885
$ws.Range("I18:I52").Value = $newRF
